$d = $word.ActiveDocument

# --- Change 1: remove the whole list-item paragraph
#     "Sistem de crafting pentru obiecte utile si decor."
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Sistem de crafting pentru obiecte utile*decor.*") {
        $p.Range.Delete() | Out-Null
        break
    }
}

# --- Change 2: "Gatit si crafting:" -> "Gatit:"
$d.Content.Find.Execute(" și crafting:", $false, $false, $false, $false, $false, $true, 1, $false, ":", 2) | Out-Null

# --- Change 3: "...speciale sau obiecte decorative." -> "...speciale."
$d.Content.Find.Execute(" sau obiecte decorative.", $false, $false, $false, $false, $false, $true, 1, $false, ".", 2) | Out-Null

# --- Change 4: wrap the standalone "a" in a grammar-proofing span
#     (splits the single " a " run into " " + <gramStart>a<gramEnd> + " ")
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*eficientiza*") {
        $xml = @'
<w:p w14:paraId="5752C0F4" w14:textId="1F5EF912" w:rsidR="002417ED" w:rsidRPr="00B245CC" w:rsidRDefault="002417ED" w:rsidP="002417ED"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3120"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">Pe </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>măsură</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ce</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>avansează</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>poate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>accesa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>echipamente</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pentru</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>eficientiza</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>procesul</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>agricol</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>
'@
        $p.Range.InsertXML($xml) | Out-Null
        break
    }
}
